# Refresh Price (D) and Volume(1h) (E) columns with the latest crypto snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.314.92"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").Value = "1.566.53"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "`'211.47"
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").Value = "`'22.25"
$ws.Range("E8").Value = "  +0.87%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").Value = "`'0.0873"
$ws.Range("E11").Value = "  +2.17%  "

$ws.Range("D12").Value = "1.788.68"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "1.561.22"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").Value = "`'3.77"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").Value = "`'0.521"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").Value = "27.314.07"
$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("D17").Value = "`'61.89"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").Value = "`'218.92"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("D19").Value = "0.0₃0708"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").Value = "`'7.47"
$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "`'4.15"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "`'9.41"
$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").Value = "`'151.45"
$ws.Range("E25").Value = "  -1.19%  "

$ws.Range("D26").Value = "`'6.64"
$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("E27").Value = "  +1.33%  "

$ws.Range("D28").Value = "`'15.04"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("D31").Value = "`'0.0473"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").Value = "1.461.25"
$ws.Range("E33").Value = "  +2.02%  "

$ws.Range("D34").Value = "`'3.18"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").Value = "`'1.12"
$ws.Range("E35").Value = "  +5.17%  "

$ws.Range("E36").Value = "  +1.21%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("D41").Value = "`'0.818"
$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("D45").Value = "`'64.52"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").Value = "1.703.39"
$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("D48").Value = "`'86.18"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").Value = "`'0.0526"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("E51").Value = "  -1.30%  "
